$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: section header ---
$ws.Range("A33").Value = " + Extra dolgok:"

# --- Values entered in the same order the author typed them, so the ---
# --- sharedStrings table indices line up with the original workbook. ---
$ws.Range("C34").Value = "43-06-67"
$ws.Range("D34").Value = "HÜVELY 20x2 P FHDH2.54-40 (T-T) 1# h=8,5mm AU"

$ws.Range("C35").Value = "43-00-12"
$ws.Range("C36").Value = "43-00-57"
$ws.Range("C37").Value = "43-05-85"

$ws.Range("D35").Value = "SZAL.KABEL 6 P ANYA FC-06 (T-T) TEHERMENTESITŐVEL LENGŐ"
$ws.Range("D36").Value = "SZAL.KABEL 10 P ANYA FC-10 (T-T) TEHERMENTESITŐVEL LENGŐ"
$ws.Range("D37").Value = "SZAL.KABEL 10 P APA 90° BHR-10 (T-T)"

$ws.Range("A34").Value = "DS"
$ws.Range("A35").Value = "DS"
$ws.Range("A36").Value = "DS"
$ws.Range("A37").Value = "DS"

$ws.Range("A38").Value = "Encoder"

$ws.Range("D38").Value = "TÁPCSATL. 2.54mm 4P ANYA HÁZ NCH254-04 (G-S)"
$ws.Range("D39").Value = "TÁPCSATL. 2.54mm 4P ANYA HÁZ NCH254-04 (G-S)"

$ws.Range("A39").Value = "Bluetooth"

$ws.Range("C38").Value = "43-09-08"
$ws.Range("C39").Value = "43-09-08"

# --- Supplier column (re-uses the existing "Lomex" shared string) ---
$ws.Range("B34").Value = "Lomex"
$ws.Range("B35").Value = "Lomex"
$ws.Range("B36").Value = "Lomex"
$ws.Range("B37").Value = "Lomex"
$ws.Range("B38").Value = "Lomex"
$ws.Range("B39").Value = "Lomex"

# --- Quantity column ---
$ws.Range("E34").Value = 2
$ws.Range("E35").Value = 8
$ws.Range("E36").Value = 8
$ws.Range("E37").Value = 8
$ws.Range("E38").Value = 2
$ws.Range("E39").Value = 2

# --- Row 43: stray formatted (empty) cell (creates the General-ish style first) ---
$ws.Range("C43").NumberFormat = "General"

# --- Reference number format for the Encoder / Bluetooth rows (text) ---
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C39").NumberFormat = "@"

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection matches the saved state ---
[void]$ws.Range("D44").Select()
